$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 17, shifting existing rows 17-42 down to 18-43
$ws.Rows.Item(17).Insert()

# Populate the new row 17 with the new weekly record
$ws.Cells.Item(17, 1).Value = 7
$ws.Cells.Item(17, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(17, 3).Value = "Ñuble"
$ws.Cells.Item(17, 4).Value = 45260
$ws.Cells.Item(17, 5).Value = 16
$ws.Cells.Item(17, 6).Value = "Fruta"
$ws.Cells.Item(17, 7).Value = 100107
$ws.Cells.Item(17, 8).Value = "Otros"
$ws.Cells.Item(17, 9).Value = 100107002
$ws.Cells.Item(17, 10).Value = "Chirimoya"
$ws.Cells.Item(17, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(17, 12).Value = "Primera"
$ws.Cells.Item(17, 13).Value = 60
$ws.Cells.Item(17, 14).Value = 19000
$ws.Cells.Item(17, 15).Value = 20000
$ws.Cells.Item(17, 16).Value = 19500
$ws.Cells.Item(17, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(17, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(17, 19).Value = 1950
$ws.Cells.Item(17, 20).Value = 10
